$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "23.495.19"
Set-TextValue "E2" "  -0.64%  "
Set-TextValue "D3" "1.648.39"
Set-TextValue "E3" "  +0.04%  "
Set-TextValue "E4" "  +0.41%  "
Set-TextValue "D5" "1.002"
Set-TextValue "E5" "  +0.36%  "
Set-TextValue "D6" "299.71"
Set-TextValue "E6" "  -1.57%  "
Set-TextValue "D7" "0.3790"
Set-TextValue "E7" "  -0.36%  "
Set-TextValue "D8" "0.3565"
Set-TextValue "E8" "  -1.17%  "
Set-TextValue "D9" "50.60"
Set-TextValue "E9" "  -2.65%  "
Set-TextValue "D10" "0.08100"
Set-TextValue "E10" "  -1.29%  "
Set-TextValue "E11" "  -1.86%  "
Set-TextValue "E12" "  +0.26%  "
Set-TextValue "D13" "22.06"
Set-TextValue "E13" "  -1.87%  "
Set-TextValue "D14" "6.415"
Set-TextValue "E14" "  -1.69%  "
Set-TextValue "D15" "7.404"
Set-TextValue "E15" "  +0.40%  "
Set-TextValue "D16" "0.00001202"
Set-TextValue "E16" "  -2.35%  "
Set-TextValue "D17" "1.655.94"
Set-TextValue "E17" "  +0.64%  "
Set-TextValue "D18" "97.22"
Set-TextValue "E18" "  +0.22%  "
Set-TextValue "D19" "0.06982"
Set-TextValue "E19" "  -0.03%  "
Set-TextValue "D20" "6.769"
Set-TextValue "E20" "  +0.50%  "
Set-TextValue "E21" "  -0.79%  "
Set-TextValue "D22" "1.001"
Set-TextValue "E22" "  +0.32%  "
Set-TextValue "D23" "12.52"
Set-TextValue "E23" "  -0.31%  "
Set-TextValue "D24" "23.527.71"
Set-TextValue "E24" "  -0.45%  "
Set-TextValue "D25" "2.488"
Set-TextValue "E25" "  -1.37%  "
Set-TextValue "D26" "2.903"
Set-TextValue "E26" "  -6.84%  "
Set-TextValue "D27" "20.96"
Set-TextValue "E27" "  -1.52%  "
Set-TextValue "D28" "152.81"
Set-TextValue "E28" "  +0.34%  "
Set-TextValue "D29" "5.219"
Set-TextValue "E29" "  +0.26%  "
Set-TextValue "D30" "133.10"
Set-TextValue "E30" "  -1.38%  "
Set-TextValue "D31" "1.839.23"
Set-TextValue "E31" "  +0.54%  "
Set-TextValue "D32" "6.951"
Set-TextValue "E32" "  +2.57%  "
Set-TextValue "D33" "2.140"
Set-TextValue "E33" "  +4.44%  "
Set-TextValue "D34" "11.89"
Set-TextValue "E34" "  +2.23%  "
Set-TextValue "D35" "1.026"
Set-TextValue "E35" "  -5.98%  "
Set-TextValue "D36" "0.02729"
Set-TextValue "E36" "  -2.67%  "
Set-TextValue "E37" "  -0.93%  "
Set-TextValue "B38" "Aptos"
Set-TextValue "C38" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D38" "13.47"
Set-TextValue "E38" "  +4.98%  "
Set-TextValue "B39" "Algorand"
Set-TextValue "C39" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D39" "0.2454"
Set-TextValue "E39" "  -2.34%  "
Set-TextValue "D40" "5.974"
Set-TextValue "E40" "  -1.93%  "
Set-TextValue "D41" "0.06845"
Set-TextValue "E41" "  -2.60%  "
Set-TextValue "D42" "0.6928"
Set-TextValue "E42" "  -1.82%  "
Set-TextValue "D43" "1.321"
Set-TextValue "E43" "  -0.77%  "
Set-TextValue "D44" "15.65"
Set-TextValue "E44" "  -1.17%  "
Set-TextValue "D45" "0.6442"
Set-TextValue "E45" "  -1.01%  "
Set-TextValue "E46" "  +0.30%  "
Set-TextValue "D47" "2.268"
Set-TextValue "E47" "  -3.03%  "
Set-TextValue "D48" "3.927"
Set-TextValue "D49" "0.07796"
Set-TextValue "E49" "  -2.29%  "
Set-TextValue "D50" "128.15"
Set-TextValue "E50" "  +0.25%  "
Set-TextValue "D51" "1.169"
Set-TextValue "E51" "  -1.73%  "
